$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 193-194 (everything currently at row 193 and below
# shifts down by two rows). Excel's default Insert behavior copies the
# formatting of the row above, which is what we want for column D's date
# style (s="2").
$ws.Rows("193:194").Insert()

# New row 193: Camote, 1a nueva(o), date 2022-10-11 (44845), origin Peru
$ws.Range("A193").Value = 11
$ws.Range("B193").Value = "Vega Monumental Concepción"
$ws.Range("C193").Value = "Bíobío"
$ws.Range("D193").Value = 44845
$ws.Range("E193").Value = 8
$ws.Range("F193").Value = 100112045
$ws.Range("G193").Value = "Zapallo"
$ws.Range("H193").Value = "Camote"
$ws.Range("I193").Value = "1a nueva(o)"
$ws.Range("J193").Value = 500
$ws.Range("K193").Value = 1000
$ws.Range("L193").Value = 1000
$ws.Range("M193").Value = 1000
$ws.Range("N193").Value = "$/kilo (volumen en unidades)"
$ws.Range("O193").Value = "Perú"
$ws.Range("P193").Value = 1000
$ws.Range("Q193").Value = 1
$ws.Range("R193").Value = "Hortaliza"

# New row 194: Camote, 2a nueva(o), date 2022-10-11 (44845), origin Peru
$ws.Range("A194").Value = 11
$ws.Range("B194").Value = "Vega Monumental Concepción"
$ws.Range("C194").Value = "Bíobío"
$ws.Range("D194").Value = 44845
$ws.Range("E194").Value = 8
$ws.Range("F194").Value = 100112045
$ws.Range("G194").Value = "Zapallo"
$ws.Range("H194").Value = "Camote"
$ws.Range("I194").Value = "2a nueva(o)"
$ws.Range("J194").Value = 300
$ws.Range("K194").Value = 800
$ws.Range("L194").Value = 800
$ws.Range("M194").Value = 800
$ws.Range("N194").Value = "$/kilo (volumen en unidades)"
$ws.Range("O194").Value = "Perú"
$ws.Range("P194").Value = 800
$ws.Range("Q194").Value = 1
$ws.Range("R194").Value = "Hortaliza"
